# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-record row is inserted above the current row 83
# (Rabanito @ Vega Central Mapocho de Santiago), shifting every
# subsequent row down by one (old row 83 -> new row 84, ..., old row
# 211 -> new row 212). The inherited data is identical to the row it
# displaces except for the Fecha (D) and Volumen (J) columns, which get
# the new reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 83; Excel shifts rows 83:211 down to 84:212.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new weekly record.
$ws.Range("A83").Value = 9
$ws.Range("B83").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C83").Value = "Metropolitana"
$ws.Range("D83").Value = 44579
$ws.Range("E83").Value = 13
$ws.Range("F83").Value = 300000001
$ws.Range("G83").Value = "Rabanito"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 6100
$ws.Range("K83").Value = 3000
$ws.Range("L83").Value = 3000
$ws.Range("M83").Value = 3000
$ws.Range("N83").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O83").Value = "Provincia de Chacabuco"
$ws.Range("P83").Value = 30
$ws.Range("Q83").Value = 100
$ws.Range("R83").Value = "Hortaliza"
